$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the simulation results table (columns B = "in pure BC", C = "in MBS", D = difference)
# to reflect the new marc_s2 simulation data.

$ws.Range("B2").Value = "+3.08 ± 0.00"
$ws.Range("C2").Value = "+3.08 ± 0.00"
$ws.Range("D2").Value = "+0.16 ± 0.00"

$ws.Range("B3").Value = "+3.08 ± 0.00"
$ws.Range("C3").Value = "+3.08 ± 0.00"

$ws.Range("D4").Value = "+0.16 ± 0.00"

$ws.Range("B5").Value = "-3.09 ± 0.00"
$ws.Range("C5").Value = "-3.09 ± 0.00"
$ws.Range("D5").Value = "-0.16 ± 0.00"

$ws.Range("C6").Value = "-0.16 ± 0.00"

$ws.Range("B7").Value = "-0.08 ± 0.00"
$ws.Range("C7").Value = "-0.00 ± 0.00"
$ws.Range("D7").Value = "-0.08 ± 0.00"

$ws.Range("B8").Value = "-0.00 ± 0.00"
$ws.Range("C8").Value = "-0.00 ± 0.00"
$ws.Range("D8").Value = "-0.00 ± 0.00"

$ws.Range("B9").Value = "-2.81 ± 0.00"
$ws.Range("C9").Value = "-2.74 ± 0.00"
$ws.Range("D9").Value = "-0.07 ± 0.00"

$ws.Range("B10").Value = "-0.19 ± 0.00"
$ws.Range("C10").Value = "-0.18 ± 0.00"
$ws.Range("D10").Value = "-0.00 ± 0.00"

$ws.Range("B11").Value = "+0.05 ± 0.00"
$ws.Range("C11").Value = "+0.05 ± 0.00"
$ws.Range("D11").Value = "+0.00 ± 0.00"

$ws.Range("B12").Value = "+5.86 ± 0.01"
$ws.Range("C12").Value = "+5.63 ± 0.01"
$ws.Range("D12").Value = "+4.42 ± 0.04"
